$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Append three new species-observation rows (62-64) to the "Artfynd" sheet,
# mirroring the layout of the existing rows above them.
# ---------------------------------------------------------------------------

# Helper: write a value as literal text, defeating Excel's automatic
# date/number recognition (used for the Startdatum/Slutdatum columns, whose
# values look like ISO dates but must stay plain text), then strip the
# temporary "Text" number format we had to apply so the cell ends up with no
# explicit style, same as every other cell in the sheet.
function Write-PlainText($rng, [string]$text) {
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

$rows = @(
    @{
        Row = 62
        A = 112358151; B = 90307; C = "Ovaliderad"; D = "VU"; E = 918
        F = "Tajgaskinn"; G = "Laurilia sulcata"; H = "(Burt) Pouzar"
        P = "Njaka, Ås lm"; Q = 523688; R = 7204521; S = 25
        T = "Västerbotten"; U = "Vilhelmina"; V = "Åsele lappmark"; W = "Vilhelmina"
        Y = "2023-09-27"; Z = "16:22"; AA = "2023-09-27"; AB = "16:22"
        AD = $false; AE = $false; AG = $false
        AW = "Yasmine Kindlund"; AX = "Yasmine Kindlund, Isak Vahlström"
    },
    @{
        Row = 63
        A = 112358152; B = 89889; C = "Ovaliderad"; D = "VU"; E = 1506
        F = "Ostticka"; G = "Skeletocutis odora"; H = "(Sacc.) Ginns"
        P = "Njaka, Ås lm"; Q = 523699; R = 7204491; S = 25
        T = "Västerbotten"; U = "Vilhelmina"; V = "Åsele lappmark"; W = "Vilhelmina"
        Y = "2023-09-27"; Z = "16:18"; AA = "2023-09-27"; AB = "16:18"
        AD = $false; AE = $false; AG = $false
        AW = "Yasmine Kindlund"; AX = "Yasmine Kindlund, Isak Vahlström"
    },
    @{
        Row = 64
        A = 112358153; B = 89881; C = "Ovaliderad"; D = "VU"; E = 2063
        F = "Grantickeporing"; G = "Skeletocutis chrysella"; H = "Niemelä"
        P = "Njaka, Ås lm"; Q = 523699; R = 7204491; S = 25
        T = "Västerbotten"; U = "Vilhelmina"; V = "Åsele lappmark"; W = "Vilhelmina"
        Y = "2023-09-27"; Z = "16:17"; AA = "2023-09-27"; AB = "16:17"
        AD = $false; AE = $false; AG = $false
        AW = "Yasmine Kindlund"; AX = "Yasmine Kindlund, Isak Vahlström"
    }
)

foreach ($r in $rows) {
    $n = $r.Row

    $ws.Range("A$n").Value = $r.A
    $ws.Range("B$n").Value = $r.B
    $ws.Range("C$n").Value = $r.C
    $ws.Range("D$n").Value = $r.D
    $ws.Range("E$n").Value = $r.E
    $ws.Range("F$n").Value = $r.F
    $ws.Range("G$n").Value = $r.G
    $ws.Range("H$n").Value = $r.H

    $ws.Range("P$n").Value = $r.P
    $ws.Range("Q$n").Value = $r.Q
    $ws.Range("R$n").Value = $r.R
    $ws.Range("S$n").Value = $r.S
    $ws.Range("T$n").Value = $r.T
    $ws.Range("U$n").Value = $r.U
    $ws.Range("V$n").Value = $r.V
    $ws.Range("W$n").Value = $r.W

    # Date/time columns stored as plain text in the source data.
    Write-PlainText $ws.Range("Y$n") $r.Y
    $ws.Range("Z$n").Value = $r.Z
    Write-PlainText $ws.Range("AA$n") $r.AA
    $ws.Range("AB$n").Value = $r.AB

    $ws.Range("AD$n").Value = $r.AD
    $ws.Range("AE$n").Value = $r.AE
    $ws.Range("AG$n").Value = $r.AG

    $ws.Range("AW$n").Value = $r.AW
    $ws.Range("AX$n").Value = $r.AX
}
